# Name changes on slides
#
# Slide 1 (title slide):
#   - Title: merge "Intermediate " + "Software Craftsmanship" into a
#     single run reading "Intermediate Software Craftsmanship".
#   - Subtitle: change presenter name/handle from
#     "Steve Smith | @ardalis" to "Jeff Valore | @CodingWithSpike".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Title 1 : "Intermediate " + "Software Craftsmanship" -> one run ----
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange

# The second run ("Software Craftsmanship") is merged away, and its text
# is appended onto the first run so the paragraph ends up as a single run.
$run2 = $titleRange.Characters(14, 22)
$run2.Text = ""

$run1 = $titleRange.Characters(1, 13)
$run1.Text = "Intermediate Software Craftsmanship"

# ---- Subtitle 2 : "Steve Smith | @ardalis" -> "Jeff Valore | @CodingWithSpike" ----
$subRange = $s.Shapes.Item(2).TextFrame.TextRange

# Work from the end of the string towards the start so earlier character
# offsets stay valid while later ones are being rewritten.

# "ardalis" (chars 16-22) -> "CodingWithSpike"
$handle = $subRange.Characters(16, 7)
$handle.Text = "CodingWithSpike"

# "Steve Smith " (chars 1-12) -> "Jeff Valore "
$name = $subRange.Characters(1, 12)
$name.Text = "Jeff Valore "

# "| " (chars 13-14) is re-assigned (same text) so it becomes its own run,
# leaving "@" (char 15) as a trailing run of its own as well.
$sep = $subRange.Characters(13, 2)
$sep.Text = "| "
